# Add the new OrangeHRM ESS test-scenario rows to the "Test Scenarios" sheet
# (requirement IDs in column B, scenario descriptions in column C, for rows
# 7 through 44 of the existing table).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Scenarios")

# Row number -> [ "Req#ID" (col B, or $null to leave blank), "Scenario" (col C) ]
$rows = [ordered]@{
    7  = @("3.1.3",  "ESS user can add photograph")
    8  = @("3.1.3",  "ESS user can view Contact Details")
    9  = @("3.1.4",  "ESS user can add Emergency Contact")
    10 = @("3.1.4",  "ESS user can add multiple Emergency Contact")
    11 = @("3.1.4",  "ESS user can delete Emergency Contact")
    12 = @("3.1.5",  "ESS user can add Dependent")
    13 = @("3.1.5",  "ESS user can add multiple Dependent")
    14 = @("3.1.5",  "ESS user can delete Dependent")
    15 = @("3.1.5",  "ESS user can add attachment to Dependent")
    16 = @("3.1.6",  "ESS user can add Immigration details")
    17 = @("3.1.6",  "ESS user can add multiple immigration details")
    18 = @("3.1.6",  "ESS user can delete immigration details")
    19 = @("3.1.6",  "ESS user can add attachment to immigration record")
    20 = @("3.1.7",  "ESS user can not change job detail ( view only )")
    21 = @("3.1.8",  "ESS user restricted to chnage Salary Components")
    22 = @("3.1.9",  "ESS user can only see the list of supervisors ")
    23 = @("3.1.10", "ESS user can add work experience")
    24 = @("3.1.10", "ESS user can add multiple work experience")
    25 = @("3.1.10", "ESS user can delete work experience")
    26 = @("3.1.10", "ESS user can add Education")
    27 = @("3.1.10", "ESS user can add multiple entries of education")
    28 = @("3.1.10", "ESS user can delete education")
    29 = @("3.1.10", "ESS user can add skill")
    30 = @("3.1.10", "ESS user can add multiple skill")
    31 = @("3.1.10", "ESS user can delete skill")
    32 = @("3.1.10", "ESS user can add language")
    33 = @("3.1.10", "ESS user can add multiple language")
    34 = @("3.1.10", "ESS user can delete language")
    35 = @("3.1.10", "ESS user can add license")
    36 = @("3.1.10", "ESS user can add multiple license")
    37 = @("3.1.10", "ESS user can delete license")
    38 = @("3.1.10", "ESS user can upload attachment")
    39 = @("3.1.10", "ESS user can upload multiple attachment")
    40 = @("3.1.10", "ESS user can delete attachment")
    41 = @("3.1.11", "ESS user can add membership")
    42 = @($null,    "ESS user can add multiple memberships")
    43 = @($null,    "ESS user can delete memberships")
    44 = @($null,    "ESS user can add attachment under memberships")
}

# Requirement-ID strings that look like dates to Excel's auto-detection
# (e.g. "3.1.10" / "3.1.11" parse as 1-Mar-2010 / 1-Mar-2011) need to be
# entered as a text formula and then converted to a plain value via
# copy / paste-special, so they land as plain text without forcing a
# number-format change on the cell (keeps styles.xml untouched).
function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

foreach ($r in $rows.Keys) {
    $pair = $rows[$r]
    $bVal = $pair[0]
    $cVal = $pair[1]

    if ($null -ne $bVal) {
        $bCell = $ws.Cells.Item($r, 2)
        if ($bVal -like "3.1.1?") {
            Set-TextValue $bCell $bVal
        } else {
            $bCell.Value = $bVal
        }
    }

    $ws.Cells.Item($r, 3).Value = $cVal
}

$excel.CutCopyMode = 0

# Match the author's final scroll position / selection on the sheet.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1
$ws.Range("C41").Select()
